$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Task 3 bullet list: three consecutive list-paragraphs ("Sending a
#    message using the xBee...", "It then deactivates the autonomous
#    behaviour...", "This then reactivates the autonomous behaviour...")
#    all become highlighted in green (wdGreen = 4). The middle paragraph
#    was previously highlighted yellow; the other two had no highlight.
# ---------------------------------------------------------------------------
$wdGreen = 4

$pSending = $d.Content
$pSending.Find.Execute("Sending a message using the", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$paraIndex1 = $d.Range($pSending.Start, $pSending.Start).Paragraphs(1).Index

$pDeactivates = $d.Content
$pDeactivates.Find.Execute("It then deactivates the", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$paraIndex2 = $d.Range($pDeactivates.Start, $pDeactivates.Start).Paragraphs(1).Index

$pReactivates = $d.Content
$pReactivates.Find.Execute("This then reactivates the", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$paraIndex3 = $d.Range($pReactivates.Start, $pReactivates.Start).Paragraphs(1).Index

$d.Paragraphs($paraIndex1).Range.HighlightColorIndex = $wdGreen
$d.Paragraphs($paraIndex2).Range.HighlightColorIndex = $wdGreen
$d.Paragraphs($paraIndex3).Range.HighlightColorIndex = $wdGreen

# ---------------------------------------------------------------------------
# 2) Split the run "an appropriate message should appear in the GUI" into
#    two runs around a relocated "_GoBack" bookmark: "...appear in th" |
#    bookmark | "e GUI".
# ---------------------------------------------------------------------------
$rSplit = $d.Content
$rSplit.Find.Execute("an appropriate message should appear in th", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPoint = $d.Range($rSplit.End, $rSplit.End)

# ---------------------------------------------------------------------------
# 3) Merge the hyperlink text ("http://andybrown.me.uk/2011/01/1" + bookmark
#    + "5" + "/the-standard-template-library-stl-for-avr-with-c-streams/")
#    back into a single run, which removes the old "_GoBack" bookmark that
#    used to sit inside it.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $hlink = $d.Hyperlinks.Item($i)
    if ($hlink.TextToDisplay -like "http://andybrown.me.uk*") {
        $hlink.TextToDisplay = "http://andybrown.me.uk/2011/01/15/the-standard-template-library-stl-for-avr-with-c-streams/"
    }
}

# Now (re)create the "_GoBack" bookmark at its new location.
$d.Bookmarks.Add("_GoBack", $splitPoint) | Out-Null
